$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto data: row, Coin, Link, Price, Volume(1h)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "31.179.70", "  +1.93%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.960.31", "  +2.25%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  +0.19%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "246.98", "  +0.80%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  +0.11%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4888", "  +1.14%  "),
    @(8, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "44.62", "  +0.27%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2970", "  +2.52%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06847", "  +0.73%  "),
    @(11, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "19.15", "  -1.94%  "),
    @(12, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "106.34", "  -5.19%  "),
    @(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07774", "  +2.58%  "),
    @(14, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.932.59", "  +0.77%  "),
    @(15, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.430", "  +0.43%  "),
    @(16, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.7108", "  +5.72%  "),
    @(17, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "286.58", "  -2.83%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "31.182.21", "  +1.96%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007771", "  +2.18%  "),
    @(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.24", "  +1.63%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.594", "  +1.32%  "),
    @(22, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.000", "  +0.05%  "),
    @(23, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.180.35", "  +0.70%  "),
    @(24, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  +0.22%  "),
    @(25, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.595", "  +2.67%  "),
    @(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.982", "  +5.23%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "168.06", "  +1.15%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.01", "  -1.38%  "),
    @(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.198", "  +5.01%  "),
    @(30, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1065", "  -0.06%  "),
    @(31, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.442", "  -0.01%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.801", "  +18.40%  "),
    @(33, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.502", "  +9.17%  "),
    @(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05027", "  +0.77%  "),
    @(35, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7670", "  +4.27%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.166", "  +1.99%  "),
    @(37, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02052", "  +1.16%  "),
    @(38, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.729", "  +0.52%  "),
    @(39, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.714", "  +1.15%  "),
    @(40, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.130", "  +5.30%  "),
    @(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.415", "  +9.77%  "),
    @(42, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8858", "  +2.39%  "),
    @(43, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "109.91", "  +0.42%  "),
    @(44, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "73.44", "  +5.75%  "),
    @(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4455", "  +0.46%  "),
    @(46, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.000", "  +0.00%  "),
    @(47, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.500", "  +3.81%  "),
    @(48, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "995.41", "  +17.69%  "),
    @(49, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1269", "  +3.55%  "),
    @(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.408", "  +2.46%  "),
    @(51, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "36.04", "  +3.75%  ")

)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
